$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @(22498.56, 0, 0, 11520, 16910.97, 0, 50929.53)
  3 = @(21985.99, 0, 30.79, 13248, 16540.65, 0, 51805.43)
  4 = @(23854.82, 2323.9, 693.4, 11952, 16854.9, 185.97, 55864.990000000005)
  5 = @(35485.13, 25914.59, 8216.36, 15552, 20551.35, 2687.62, 108407.04999999999)
  6 = @(87036.77, 135940.31, 42676.89, 50688, 48075.07, 23981.08, 388398.12000000005)
  7 = @(174123.43, 294064.14, 91546.6, 165312, 140757.72, 119035.97, 984839.86)
  8 = @(309196.44, 637007.94, 169399.69, 319968, 275531.02, 289989.18, 2001092.2699999998)
  9 = @(342077.9, 919414, 219039.01, 236160, 202643.26, 175946.06, 2095280.23)
  10 = @(91173.6, 213868.43, 55604.07, 66384, 51454.45, 41025.91, 519510.4600000001)
  11 = @(45477.94, 83172.57, 23725.63, 42768, 29809.51, 25275.5, 250229.15000000002)
  12 = @(35290.08, 58541.06, 19032.29, 29376, 21597.93, 15357.05, 179194.40999999997)
  13 = @(51451.85, 98694.22, 34182.15, 28656, 19262.13, 8777.88, 241024.23)
  14 = @(97605.65, 172450.76, 79397.75, 30240, 25785.34, 4772.27, 410251.7700000001)
  15 = @(138452.33, 152199.9, 109267.26, 44352, 46375.51, 2412.73, 493059.73)
  16 = @(169718.98, 165781.06, 104992.97, 82512, 87610, 906.99, 611522)
  17 = @(225856.51, 228133.64, 116714.5, 124848, 148030.49, 447.45, 844030.59)
  18 = @(272785.97, 292342.57, 105406.58, 181296, 219152.19, 83.2, 1071066.51)
  19 = @(304320.24, 345825.14, 80827.76, 233568, 279286.04, 43.95, 1243871.13)
  20 = @(347711.62, 584450.36, 97618.75, 260640, 315581.01, 2.84, 1606004.58)
  21 = @(336031.42, 476414.17, 84929.69, 259344, 310325.02, 0, 1467044.3)
  22 = @(291007.08, 328205.12, 79629.07, 215136, 272284.27, 6.25, 1186267.79)
  23 = @(248119.2, 271564.83, 106462.05, 175968, 206664.28, 0, 1008778.3600000001)
  24 = @(174486.31, 181982.86, 88774.3, 112032, 125266.81, 0, 682542.28)
  25 = @(112628.88, 111315.81, 66530.67, 62784, 62623.86, 0, 415883.22)
  26 = @(3958376.7000000007, 5779607.38, 1784698.23, 2774304, 2958973.78, 710937.8999999999, 17966897.99)
}

foreach ($r in $data.Keys) {
  $vals = $data[$r]
  for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item([int]$r, $i + 2).Value = $vals[$i]
  }
}
